$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the four values B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: D2 value moved to B2 (new value), D2 cleared
$ws.Range("B2").Value = 42.093751076109172
$ws.Range("D2").ClearContents()

# Row 3: B3 cleared, C3 value changed
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 41.069221641349202

# Update the selected range to match the new selection B1:E3
$ws.Range("B1:E3").Select()
